# Regenerate save_data to use K instead of Strike#, regen std/mean, calc and write s_vals.
# This updates column G ("K") values for rows 2-40 on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 1
    3  = 1
    4  = 2
    5  = 1
    6  = 1
    7  = 3
    8  = 1
    9  = 1
    10 = 2
    11 = 1
    12 = 1
    13 = 1
    14 = 3
    15 = 3
    16 = 0
    17 = 1
    18 = 1
    19 = 0
    20 = 2
    21 = 0
    22 = 1
    23 = 2
    24 = 0
    25 = 1
    26 = 0
    27 = 1
    28 = 1
    29 = 0
    30 = 2
    31 = 5
    32 = 3
    33 = 2
    34 = 3
    35 = 1
    36 = 2
    37 = 1
    38 = 2
    39 = 1
    40 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
